# Bugfix DGV1 line 8 is pipe diameter NOT montage
# Add a note under the header explaining the broken screw, and
# update the visible selection to cover the header + new note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Broken screw after 9 years"

$ws.Range("A1:D20").Select()
